$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 281 (existing rows 281-296 shift
# down to 283-298, preserving their data/formatting untouched).
$ws.Rows("281:282").Insert()

# New row 281: Primera quality, week of 2022-07-11 (serial 44753)
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 44753
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100112037
$ws.Range("G281").Value = "Cebollín"
$ws.Range("H281").Value = "Sin especificar"
$ws.Range("I281").Value = "Primera"
$ws.Range("J281").Value = 35
$ws.Range("K281").Value = 12000
$ws.Range("L281").Value = 12000
$ws.Range("M281").Value = 12000
$ws.Range("N281").Value = "$/paquete 36 unidades"
$ws.Range("O281").Value = "Región Metropolitana"
$ws.Range("P281").Value = 333
$ws.Range("Q281").Value = 36
$ws.Range("R281").Value = "Hortaliza"

# New row 282: Segunda quality, same date
$ws.Range("A282").Value = 4
$ws.Range("B282").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C282").Value = "Los Lagos"
$ws.Range("D282").Value = 44753
$ws.Range("E282").Value = 10
$ws.Range("F282").Value = 100112037
$ws.Range("G282").Value = "Cebollín"
$ws.Range("H282").Value = "Sin especificar"
$ws.Range("I282").Value = "Segunda"
$ws.Range("J282").Value = 35
$ws.Range("K282").Value = 9000
$ws.Range("L282").Value = 9000
$ws.Range("M282").Value = 9000
$ws.Range("N282").Value = "$/paquete 36 unidades"
$ws.Range("O282").Value = "Región Metropolitana"
$ws.Range("P282").Value = 250
$ws.Range("Q282").Value = 36
$ws.Range("R282").Value = "Hortaliza"
